$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F data (rows 25-29), mirroring existing B/C/D/E structure
$ws.Range("F25").Value = 1000
$ws.Range("F26").Value = "XGB"
$ws.Range("F27").Value = 30
$ws.Range("F28").Value = 0.597
$ws.Range("F29").Value = 0.608

# Update selection to match target state
$ws.Range("P15").Select() | Out-Null
